# feat: add NZZ article line boundaries
#
# Fills in the "first_row" (G) / "last_row" (H) line-boundary columns for
# rows 2-58 (these were already present for rows 59-189), adds a remark
# in I47, and updates the sheet view (freeze the first column, scroll down,
# move the active selection) to match the author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. First_row (G2:G58) / last_row (H2:H58) values
# ---------------------------------------------------------------------
$gValues = @(4,4,4,4,4,5,5,7,4,3,4,4,2,4,2,2,3,2,2,3,3,2,4,4,4,2,3,3,3,2,2,3,4,4,5,4,2,1,2,3,3,2,4,2,3,2,2,3,3,3,2,3,5,2,3,4,2)
$hValues = @(9,7,9,16,9,8,14,19,19,32,21,10,12,8,5,5,7,6,4,9,7,5,7,6,6,2,3,10,4,2,22,7,7,8,9,7,2,3,2,7,8,12,15,5,4,6,5,7,8,5,3,5,26,4,6,6,7)

$rowCount = $gValues.Length

$gArr = New-Object 'object[,]' $rowCount,1
$hArr = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $gArr[$i,0] = $gValues[$i]
    $hArr[$i,0] = $hValues[$i]
}

$ws.Range("G2:G58").Value = $gArr
$ws.Range("H2:H58").Value = $hArr

# ---------------------------------------------------------------------
# 2. Remark for row 47 ("xxx falsch, vertauscht?")
# ---------------------------------------------------------------------
$ws.Range("I47").Value = "xxx falsch, vertauscht?"

# ---------------------------------------------------------------------
# 3. Sheet view: freeze first column, scroll to row 54, select G57
# ---------------------------------------------------------------------
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G57").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 54
